# Refresh the run timestamps recorded in column Z ("timestamp") of the
# Log_Muestras sheet. Each batch of rows that shares the same original
# microsecond-precision timestamp is re-stamped with the timestamp of the
# newer run (commit: "mas cambios sobre lanotebook").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z45").Value   = "2025-10-17T07:09:41.690343"
$ws.Range("Z46:Z74").Value  = "2025-10-17T07:09:41.790251"
$ws.Range("Z75:Z102").Value = "2025-10-17T07:09:41.879860"
$ws.Range("Z103:Z106").Value = "2025-10-17T07:09:41.971287"
$ws.Range("Z107:Z109").Value = "2025-10-17T07:09:41.972287"
$ws.Range("Z110:Z112").Value = "2025-10-17T07:09:41.973294"
